$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 and Row 3 effectively swap their Fecha, Volumen, Precio minimo,
# Precio maximo, Precio promedio ponderado and Precio $/Kg values.

# Row 2 (new values)
$ws.Range("D2").Value = 44209
$ws.Range("M2").Value = 100
$ws.Range("N2").Value = 10000
$ws.Range("O2").Value = 11000
$ws.Range("P2").Value = 10500
$ws.Range("S2").Value = 750

# Row 3 (new values)
$ws.Range("D3").Value = 44217
$ws.Range("M3").Value = 200
$ws.Range("N3").Value = 11000
$ws.Range("O3").Value = 12000
$ws.Range("P3").Value = 11500
$ws.Range("S3").Value = 821
